$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 157.11111
$ws.Range("I11").Value = 157.11111
$ws.Range("K11").Value = 157.11111
$ws.Range("M11").Value = -17.11111

$ws.Range("H28").Value = 626.6667
$ws.Range("I28").Value = 455
$ws.Range("K28").Value = 455
$ws.Range("M28").Value = 30

$ws.Range("H29").Value = 2656.5715
$ws.Range("J29").Value = 2932.8333
$ws.Range("L29").Value = 8798.499899999999
$ws.Range("N29").Value = -9360.499899999999

$ws.Range("H32").Value = 2858.8
$ws.Range("I32").Value = 2858.8
$ws.Range("K32").Value = 2858.8
$ws.Range("M32").Value = -2532.8

$ws.Range("H116").Value = 6514.5713
$ws.Range("I116").Value = 5562.1665
$ws.Range("J116").Value = 7228.875
$ws.Range("K116").Value = 5562.1665
$ws.Range("L116").Value = 7228.875
$ws.Range("M116").Value = -2120.1665
$ws.Range("N116").Value = -14112.875

$ws.Range("H118").Value = 207.5
$ws.Range("I118").Value = 207.5
$ws.Range("K118").Value = 622.5
$ws.Range("M118").Value = 1034.5

$ws.Range("H129").Value = 20715.2
$ws.Range("I129").Value = 450.2857
$ws.Range("K129").Value = 1350.8571
$ws.Range("M129").Value = 3649.1429

$ws.Range("H132").Value = 1422
$ws.Range("I132").Value = 1422
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4266
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1736
$ws.Range("N132").ClearContents()

$ws.Range("H138").Value = 4348.467
$ws.Range("J138").Value = 4401.841
$ws.Range("L138").Value = 13205.523
$ws.Range("N138").Value = -23485.523

$ws.Range("H141").Value = 1786.4375
$ws.Range("I141").Value = 1105.5333
$ws.Range("K141").Value = 3316.5999
$ws.Range("M141").Value = 1863.4001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10140.667
$ws.Range("I32").Value = 7854.615
$ws.Range("K32").Value = 7854.615
$ws.Range("M32").Value = -7567.615

$ws.Range("H45").Value = 2030.5
$ws.Range("I45").Value = 1977.6
$ws.Range("K45").Value = 1977.6
$ws.Range("M45").Value = -1600.6

$ws.Range("H61").Value = 1428.5
$ws.Range("I61").Value = 1391.3077
$ws.Range("K61").Value = 1391.3077
$ws.Range("M61").Value = -1179.3077

$ws.Range("H132").Value = 2090.875
$ws.Range("I132").Value = 2055.0967
$ws.Range("J132").Value = 3200
$ws.Range("K132").Value = 6165.2901
$ws.Range("L132").Value = 9600
$ws.Range("M132").Value = -3635.2901
$ws.Range("N132").Value = -14660

$ws.Range("H136").Value = 1428.5
$ws.Range("I136").Value = 1391.3077
$ws.Range("K136").Value = 4173.9231
$ws.Range("M136").Value = -1623.9231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 83333336
$ws.Range("I7").Value = 83333336
$ws.Range("K7").Value = 83333336
$ws.Range("M7").Value = -83333223

$ws.Range("H20").Value = 1148.7142
$ws.Range("I20").Value = 1140.1666
$ws.Range("J20").Value = 1200
$ws.Range("K20").Value = 1140.1666
$ws.Range("L20").Value = 1200
$ws.Range("M20").Value = -893.1666
$ws.Range("N20").Value = -1694

$ws.Range("H86").Value = 1729.2858
$ws.Range("I86").Value = 1684.1666
$ws.Range("K86").Value = 1684.1666
$ws.Range("M86").Value = -561.1666

$ws.Range("H89").Value = 1729.2858
$ws.Range("I89").Value = 1684.1666
$ws.Range("K89").Value = 8420.833000000001
$ws.Range("M89").Value = -2804.833000000001

$ws.Range("H94").Value = 3670
$ws.Range("I94").Value = 2500
$ws.Range("J94").Value = 6010
$ws.Range("K94").Value = 2500
$ws.Range("L94").Value = 6010
$ws.Range("M94").Value = -2049
$ws.Range("N94").Value = -6912

$ws.Range("H99").Value = 5066.6665
$ws.Range("I99").Value = 4880
$ws.Range("K99").Value = 4880
$ws.Range("M99").Value = -3382

$ws.Range("H100").Value = 22870.666
$ws.Range("J100").Value = 22870.666
$ws.Range("L100").Value = 22870.666
$ws.Range("N100").Value = -25034.666

$ws.Range("H134").Value = 1792.091
$ws.Range("I134").Value = 1471.4
$ws.Range("K134").Value = 4414.200000000001
$ws.Range("M134").Value = -1879.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6622.909
$ws.Range("J31").Value = 8268
$ws.Range("L31").Value = 8268
$ws.Range("N31").Value = -8858

$ws.Range("H34").Value = 6622.909
$ws.Range("J34").Value = 8268
$ws.Range("L34").Value = 8268
$ws.Range("N34").Value = -8672

$ws.Range("H35").Value = 400989.94
$ws.Range("I35").Value = 471573.53
$ws.Range("J35").Value = 1016.3333
$ws.Range("K35").Value = 471573.53
$ws.Range("L35").Value = 1016.3333
$ws.Range("M35").Value = -471279.53
$ws.Range("N35").Value = -1604.3333

$ws.Range("H88").Value = 99500
$ws.Range("I88").Value = 99500
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 99500
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("M88").Value = -99094

$ws.Range("H91").Value = 99500
$ws.Range("I91").Value = 99500
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 99500
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("M91").Value = -98096

$ws.Range("H99").Value = 10163.75
$ws.Range("I99").Value = 7306.3335
$ws.Range("K99").Value = 7306.3335
$ws.Range("M99").Value = -5808.3335

$ws.Range("H105").Value = 1445.5834
$ws.Range("I105").Value = 980.875
$ws.Range("J105").Value = 2375
$ws.Range("K105").Value = 980.875
$ws.Range("L105").Value = 2375
$ws.Range("M105").Value = 766.125
$ws.Range("N105").Value = -5869

$ws.Range("H126").Value = 10163.75
$ws.Range("I126").Value = 7306.3335
$ws.Range("K126").Value = 21919.0005
$ws.Range("M126").Value = -19449.0005

$ws.Range("H132").Value = 3737.923
$ws.Range("I132").Value = 2037.5714
$ws.Range("K132").Value = 6112.7142
$ws.Range("M132").Value = -3582.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 129.66667
$ws.Range("I8").Value = 129.66667
$ws.Range("K8").Value = 389.00001
$ws.Range("M8").Value = -250.00001

$ws.Range("H12").Value = 18.714285
$ws.Range("J12").Value = 16.4
$ws.Range("L12").Value = 49.2
$ws.Range("N12").Value = -395.2

$ws.Range("H92").Value = 687
$ws.Range("I92").Value = 536.5
$ws.Range("J92").Value = 837.5
$ws.Range("K92").Value = 1609.5
$ws.Range("L92").Value = 2512.5
$ws.Range("M92").Value = -361.5
$ws.Range("N92").Value = -5008.5

$ws.Range("H107").Value = 672.36365
$ws.Range("J107").Value = 652.34375
$ws.Range("L107").Value = 1957.03125
$ws.Range("N107").Value = -5797.03125

$ws.Range("H132").Value = 2041.5555
$ws.Range("I132").Value = 1323.5
$ws.Range("K132").Value = 11911.5
$ws.Range("M132").Value = -9381.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 885.0714
$ws.Range("I107").Value = 941.2
$ws.Range("J107").Value = 820.3077
$ws.Range("K107").Value = 941.2
$ws.Range("L107").Value = 820.3077
$ws.Range("M107").Value = 978.8
$ws.Range("N107").Value = -4660.3077

$ws.Range("H113").Value = 3479.25
$ws.Range("I113").Value = 2484.875
$ws.Range("J113").Value = 4473.625
$ws.Range("K113").Value = 2484.875
$ws.Range("L113").Value = 4473.625
$ws.Range("M113").Value = -314.875
$ws.Range("N113").Value = -8813.625

$ws.Range("H132").Value = 1626.7858
$ws.Range("I132").Value = 1371.6154
$ws.Range("K132").Value = 4114.8462
$ws.Range("M132").Value = -1584.8462

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1869
$ws.Range("I7").Value = 2016.6
$ws.Range("K7").Value = 2016.6
$ws.Range("M7").Value = -1904.6

$ws.Range("H22").Value = 3823.75
$ws.Range("I22").Value = 3836.125
$ws.Range("J22").Value = 3799
$ws.Range("K22").Value = 3836.125
$ws.Range("L22").Value = 3799
$ws.Range("M22").Value = -3541.125
$ws.Range("N22").Value = -4389

$ws.Range("H27").Value = 3823.75
$ws.Range("I27").Value = 3836.125
$ws.Range("J27").Value = 3799
$ws.Range("K27").Value = 3836.125
$ws.Range("L27").Value = 3799
$ws.Range("M27").Value = -3729.125
$ws.Range("N27").Value = -4013

$ws.Range("H126").Value = 1869
$ws.Range("I126").Value = 2016.6
$ws.Range("K126").Value = 6049.799999999999
$ws.Range("M126").Value = -3579.799999999999

$ws.Range("H132").Value = 3491.7896
$ws.Range("I132").Value = 2905.138
$ws.Range("J132").Value = 5382.1113
$ws.Range("K132").Value = 8715.414000000001
$ws.Range("L132").Value = 16146.3339
$ws.Range("M132").Value = -6185.414000000001
$ws.Range("N132").Value = -21206.3339

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2011.2632
$ws.Range("I136").Value = 724.16
$ws.Range("J136").Value = 4486.4614
$ws.Range("K136").Value = 2172.48
$ws.Range("L136").Value = 13459.3842
$ws.Range("M136").Value = 377.52
$ws.Range("N136").Value = -18559.3842
